# Append 25 more days (2021-05-17 .. 2021-06-10) of scheduled/tracked flight
# counts to the "Ark1" table, extending it from row 406 to row 431 (columns
# A:D): DateTime (text), Scheduled flights, Tracked flights, and the
# Tracked/Scheduled ratio formula.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Copy formatting (styles/number formats) from the last existing data block
# down into the new rows so the appended rows inherit the same look as
# the rest of the table (text-formatted date column, number columns, % column).
$ws.Range("A382:D406").Copy($ws.Range("A407:D431"))

$ws.Range("A407").Value = "2021-05-17"
$ws.Range("B407").Value = 59
$ws.Range("C407").Value = 55
$ws.Range("D407").Formula = "=C407/B407"
$ws.Range("A408").Value = "2021-05-18"
$ws.Range("B408").Value = 59
$ws.Range("C408").Value = 55
$ws.Range("D408").Formula = "=C408/B408"
$ws.Range("A409").Value = "2021-05-19"
$ws.Range("B409").Value = 56
$ws.Range("C409").Value = 53
$ws.Range("D409").Formula = "=C409/B409"
$ws.Range("A410").Value = "2021-05-20"
$ws.Range("B410").Value = 70
$ws.Range("C410").Value = 68
$ws.Range("D410").Formula = "=C410/B410"
$ws.Range("A411").Value = "2021-05-21"
$ws.Range("B411").Value = 69
$ws.Range("C411").Value = 65
$ws.Range("D411").Formula = "=C411/B411"
$ws.Range("A412").Value = "2021-05-22"
$ws.Range("B412").Value = 60
$ws.Range("C412").Value = 57
$ws.Range("D412").Formula = "=C412/B412"
$ws.Range("A413").Value = "2021-05-23"
$ws.Range("B413").Value = 62
$ws.Range("C413").Value = 57
$ws.Range("D413").Formula = "=C413/B413"
$ws.Range("A414").Value = "2021-05-24"
$ws.Range("B414").Value = 56
$ws.Range("C414").Value = 53
$ws.Range("D414").Formula = "=C414/B414"
$ws.Range("A415").Value = "2021-05-25"
$ws.Range("B415").Value = 64
$ws.Range("C415").Value = 61
$ws.Range("D415").Formula = "=C415/B415"
$ws.Range("A416").Value = "2021-05-26"
$ws.Range("B416").Value = 62
$ws.Range("C416").Value = 61
$ws.Range("D416").Formula = "=C416/B416"
$ws.Range("A417").Value = "2021-05-27"
$ws.Range("B417").Value = 66
$ws.Range("C417").Value = 65
$ws.Range("D417").Formula = "=C417/B417"
$ws.Range("A418").Value = "2021-05-28"
$ws.Range("B418").Value = 68
$ws.Range("C418").Value = 65
$ws.Range("D418").Formula = "=C418/B418"
$ws.Range("A419").Value = "2021-05-29"
$ws.Range("B419").Value = 48
$ws.Range("C419").Value = 46
$ws.Range("D419").Formula = "=C419/B419"
$ws.Range("A420").Value = "2021-05-30"
$ws.Range("B420").Value = 41
$ws.Range("C420").Value = 41
$ws.Range("D420").Formula = "=C420/B420"
$ws.Range("A421").Value = "2021-05-31"
$ws.Range("B421").Value = 63
$ws.Range("C421").Value = 60
$ws.Range("D421").Formula = "=C421/B421"
$ws.Range("A422").Value = "2021-06-01"
$ws.Range("B422").Value = 66
$ws.Range("C422").Value = 60
$ws.Range("D422").Formula = "=C422/B422"
$ws.Range("A423").Value = "2021-06-02"
$ws.Range("B423").Value = 62
$ws.Range("C423").Value = 59
$ws.Range("D423").Formula = "=C423/B423"
$ws.Range("A424").Value = "2021-06-03"
$ws.Range("B424").Value = 78
$ws.Range("C424").Value = 70
$ws.Range("D424").Formula = "=C424/B424"
$ws.Range("A425").Value = "2021-06-04"
$ws.Range("B425").Value = 77
$ws.Range("C425").Value = 72
$ws.Range("D425").Formula = "=C425/B425"
$ws.Range("A426").Value = "2021-06-05"
$ws.Range("B426").Value = 58
$ws.Range("C426").Value = 55
$ws.Range("D426").Formula = "=C426/B426"
$ws.Range("A427").Value = "2021-06-06"
$ws.Range("B427").Value = 69
$ws.Range("C427").Value = 68
$ws.Range("D427").Formula = "=C427/B427"
$ws.Range("A428").Value = "2021-06-07"
$ws.Range("B428").Value = 67
$ws.Range("C428").Value = 63
$ws.Range("D428").Formula = "=C428/B428"
$ws.Range("A429").Value = "2021-06-08"
$ws.Range("B429").Value = 68
$ws.Range("C429").Value = 66
$ws.Range("D429").Formula = "=C429/B429"
$ws.Range("A430").Value = "2021-06-09"
$ws.Range("B430").Value = 74
$ws.Range("C430").Value = 69
$ws.Range("D430").Formula = "=C430/B430"
$ws.Range("A431").Value = "2021-06-10"
$ws.Range("B431").Value = 81
$ws.Range("C431").Value = 80
$ws.Range("D431").Formula = "=C431/B431"

# Update the visible selection to match the edited workbook state.
$ws.Range("G430").Select()
